$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the (erroneous) "Client" values that were placed in column G
# for rows 7, 9, 11, 15 and 19 - these cells are fully cleared (value +
# style), not just blanked, matching how the source file no longer
# contains these <c> elements at all.
$ws.Range("G7").Clear()
$ws.Range("G9").Clear()
$ws.Range("G11").Clear()
$ws.Range("G15").Clear()
$ws.Range("G19").Clear()

# Update the active selection left over from editing, now resting on G7.
$ws.Range("G7").Select()
